$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing "Datum/Zeit/Tätigkeit" table from columns D:F to G:J, ---
# --- adding a new "Einheit" column (I) and a new "genau" column (K).        ---

# Header row (row 5): copy the old D5/E5/F5 cells (keeps the same shared
# strings / style) into their new home, then fill the two new header cells.
$ws.Range("D5").Copy($ws.Range("G5"))
$ws.Range("E5").Copy($ws.Range("H5"))
$ws.Range("F5").Copy($ws.Range("J5"))
$ws.Range("I5").Value = "Einheit"
$ws.Range("K5").Value = "genau"

# Data rows: re-enter the dates (via Copy so the existing date style is
# reused instead of minting a new number format), hours as plain numbers,
# "Stunden" as the unit, and the task descriptions.
$ws.Range("D6").Copy($ws.Range("G6"))
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = "Stunden"
$ws.Range("J6").Value = "Recherche, Konzept Erstellung"

$ws.Range("D7").Copy($ws.Range("G7"))
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = "Stunden"
$ws.Range("J7").Value = "Konzept Weiterarbeit"

$ws.Range("D8").Copy($ws.Range("G8"))
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = "Stunden"
$ws.Range("J8").Value = "GitHub Repo erstellt"

$ws.Range("D9").Copy($ws.Range("G9"))
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = "Stunden"
$ws.Range("J9").Value = "Code-Basis erstellt"

# Clear the old D:F table now that its contents live in G:J.
$ws.Range("D5:F9").Clear()

# New rows of time-tracking entries.
$ws.Range("G10").Value = "2/13/2019"
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = "Stunden"
$ws.Range("J10").Value = "Änderungen am Code, Namen"

$ws.Range("G11").Value = "6/18/2019"
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = "Stunden"
$ws.Range("J11").Value = "Erstellen der Dokumentation"

$ws.Range("G12").Value = "6/25/2019"
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = "Stunden"
$ws.Range("J12").Value = "Anlegen des Projektmanagements in Trello"

$ws.Range("G13").Value = "6/26/2019"
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = "Stunden"
$ws.Range("J13").Value = "Diplomarbeitsbesprechung"

$ws.Range("G14").Value = "7/1/2019"
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = "Stunden"
$ws.Range("J14").Value = "Statusbericht"
$ws.Range("K14").Value = "Statusbericht geschrieben"

$ws.Range("G15").Value = "7/8/2019"
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = "Stunden"
$ws.Range("J15").Value = "Dokumentation"

$ws.Range("G16").Value = "7/14/2019"
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = "Stunden"
$ws.Range("J16").Value = "Programmieren"
$ws.Range("K16").Value = "Erste Versuche in der IDE mit fremden Librarys"

$ws.Range("G17").Value = "7/15/2019"
$ws.Range("H17").Value = 0.5
$ws.Range("I17").Value = "Stunden"
$ws.Range("J17").Value = "Dokumentation"
$ws.Range("K17").Value = "Bt Bibliothek über Maven hinzugefügt"

# Re-apply the short-date style (the "G10:G17" cells above were entered as
# text dates so Excel would parse them; make sure they all match the
# existing date format used by G6:G9) by copying the existing date style.
$ws.Range("G6").Copy($ws.Range("G10:G17"))
$ws.Range("G10").Value = "2/13/2019"
$ws.Range("G11").Value = "6/18/2019"
$ws.Range("G12").Value = "6/25/2019"
$ws.Range("G13").Value = "6/26/2019"
$ws.Range("G14").Value = "7/1/2019"
$ws.Range("G15").Value = "7/8/2019"
$ws.Range("G16").Value = "7/14/2019"
$ws.Range("G17").Value = "7/15/2019"

# New summary block in columns A:C.
$ws.Range("A6").Value = "Summe der Stunden:"
$ws.Range("B6").Formula = "=SUM(H:H)"
$ws.Range("C6").Value = "Stunden"

$ws.Range("A7").Value = "Arbeitszeit pro Tag"
$ws.Range("B7").Formula = "=((180 -B6)/(DAYS(DATE(2020,4,3), TODAY())))*60"
$ws.Range("C7").Value = "Stunden"

# Column widths: A is new, G/J reuse the widths the old D/F columns had.
$ws.Columns.Item(1).ColumnWidth = 22

# Selection, to match the recorded cursor position after the edit.
$ws.Range("K15").Select()
